$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 5720.6665
$ws.Range("I82").Value = 1082
$ws.Range("K82").Value = 3246
$ws.Range("M82").Value = -2840

$ws.Range("H85").Value = 5720.6665
$ws.Range("I85").Value = 1082
$ws.Range("K85").Value = 3246
$ws.Range("M85").Value = -1842

$ws.Range("H113").Value = 1097
$ws.Range("J113").Value = 1097
$ws.Range("L113").Value = 1097
$ws.Range("N113").Value = -7605

$ws.Range("H132").Value = 7784.048
$ws.Range("I132").Value = 7784.048
$ws.Range("K132").Value = 23352.144
$ws.Range("M132").Value = -20822.144

$ws.Range("H134").Value = 68199.56
$ws.Range("J134").Value = 68199.56
$ws.Range("L134").Value = 68199.56
$ws.Range("N134").Value = -78339.56

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 761.9286
$ws.Range("I2").Value = 761.087
$ws.Range("J2").Value = 765.8
$ws.Range("K2").Value = 761.087
$ws.Range("L2").Value = 765.8
$ws.Range("M2").Value = -648.087
$ws.Range("N2").Value = -991.8

$ws.Range("H32").Value = 3829.4368
$ws.Range("I32").Value = 2748.9136
$ws.Range("K32").Value = 2748.9136
$ws.Range("M32").Value = -2461.9136

$ws.Range("H61").Value = 4610.207
$ws.Range("I61").Value = 1310.6471
$ws.Range("K61").Value = 1310.6471
$ws.Range("M61").Value = -1098.6471

$ws.Range("H102").Value = 3726.2122
$ws.Range("I102").Value = 3475.7778
$ws.Range("J102").Value = 4853.1665
$ws.Range("K102").Value = 3475.7778
$ws.Range("L102").Value = 4853.1665
$ws.Range("M102").Value = -1853.7778
$ws.Range("N102").Value = -8097.1665

$ws.Range("H116").Value = 761.9286
$ws.Range("I116").Value = 761.087
$ws.Range("J116").Value = 765.8
$ws.Range("K116").Value = 761.087
$ws.Range("L116").Value = 765.8
$ws.Range("M116").Value = 1532.913
$ws.Range("N116").Value = -5353.8

$ws.Range("H132").Value = 2416.9211
$ws.Range("I132").Value = 1998.1
$ws.Range("J132").Value = 3987.5
$ws.Range("K132").Value = 5994.299999999999
$ws.Range("L132").Value = 11962.5
$ws.Range("M132").Value = -3464.299999999999
$ws.Range("N132").Value = -17022.5

$ws.Range("H136").Value = 4610.207
$ws.Range("I136").Value = 1310.6471
$ws.Range("K136").Value = 3931.9413
$ws.Range("M136").Value = -1381.9413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 761.9286
$ws.Range("I3").Value = 761.087
$ws.Range("J3").Value = 765.8
$ws.Range("K3").Value = 761.087
$ws.Range("L3").Value = 765.8
$ws.Range("M3").Value = -647.087
$ws.Range("N3").Value = -993.8

$ws.Range("H20").Value = 28739416
$ws.Range("I20").Value = 36235684
$ws.Range("K20").Value = 36235684
$ws.Range("M20").Value = -36235437

$ws.Range("H134").Value = 2719.3171
$ws.Range("I134").Value = 2458.9333
$ws.Range("K134").Value = 7376.7999
$ws.Range("M134").Value = -4841.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3792.3147
$ws.Range("I31").Value = 2952.7856
$ws.Range("K31").Value = 2952.7856
$ws.Range("M31").Value = -2657.7856

$ws.Range("H34").Value = 3792.3147
$ws.Range("I34").Value = 2952.7856
$ws.Range("K34").Value = 2952.7856
$ws.Range("M34").Value = -2750.7856

$ws.Range("H107").Value = 525
$ws.Range("I107").Value = 534.2917
$ws.Range("K107").Value = 534.2917
$ws.Range("M107").Value = 1385.7083

$ws.Range("H134").Value = 2063.7354
$ws.Range("I134").Value = 1985
$ws.Range("J134").Value = 2208.0833
$ws.Range("K134").Value = 5955
$ws.Range("L134").Value = 6624.249899999999
$ws.Range("M134").Value = -3420
$ws.Range("N134").Value = -11694.2499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1399.8
$ws.Range("I41").Value = 499.75
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 1499.25
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = -1161.25
$ws.Range("N41").Value = -15676

$ws.Range("H68").Value = 12504614
$ws.Range("J68").Value = 25001556
$ws.Range("L68").Value = 75004668
$ws.Range("N68").Value = -75006290

$ws.Range("H71").Value = 12504614
$ws.Range("J71").Value = 25001556
$ws.Range("L71").Value = 225014004
$ws.Range("N71").Value = -225022116

$ws.Range("H80").Value = 7192.7144
$ws.Range("I80").Value = 7483
$ws.Range("J80").Value = 6975
$ws.Range("K80").Value = 22449
$ws.Range("L80").Value = 20925
$ws.Range("M80").Value = -21513
$ws.Range("N80").Value = -22797

$ws.Range("H83").Value = 7192.7144
$ws.Range("I83").Value = 7483
$ws.Range("J83").Value = 6975
$ws.Range("K83").Value = 67347
$ws.Range("L83").Value = 62775
$ws.Range("M83").Value = -62667
$ws.Range("N83").Value = -72135

$ws.Range("H124").Value = 6257.25
$ws.Range("I124").Value = 6257.25
$ws.Range("K124").Value = 18771.75
$ws.Range("M124").Value = -13861.75

$ws.Range("H125").Value = 4000
$ws.Range("I125").Value = 3000
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 15000
$ws.Range("M125").Value = -4080
$ws.Range("N125").Value = -24840

$ws.Range("H131").Value = 5137.884
$ws.Range("J131").Value = 2043
$ws.Range("L131").Value = 6129
$ws.Range("N131").Value = -16209

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5389.9565
$ws.Range("I97").Value = 2997.0557
$ws.Range("K97").Value = 2997.0557
$ws.Range("M97").Value = -2501.0557

$ws.Range("H132").Value = 2259.2031
$ws.Range("I132").Value = 1886.875
$ws.Range("J132").Value = 3376.1875
$ws.Range("K132").Value = 5660.625
$ws.Range("L132").Value = 10128.5625
$ws.Range("M132").Value = -3130.625
$ws.Range("N132").Value = -15188.5625

$ws.Range("H136").Value = 11806.294
$ws.Range("J136").Value = 11806.294
$ws.Range("L136").Value = 35418.882
$ws.Range("N136").Value = -40518.882

$ws.Range("H139").Value = 99940.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 99940.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 99940.75
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -110220.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 707.4737
$ws.Range("I16").Value = 631.2353000000001
$ws.Range("K16").Value = 631.2353000000001
$ws.Range("M16").Value = -461.2353000000001

$ws.Range("H55").Value = 301
$ws.Range("I55").Value = 297.92307
$ws.Range("J55").Value = 307.66666
$ws.Range("K55").Value = 297.92307
$ws.Range("L55").Value = 307.66666
$ws.Range("M55").Value = -124.92307
$ws.Range("N55").Value = -653.66666

$ws.Range("H132").Value = 5030.2144
$ws.Range("I132").Value = 2944.3809
$ws.Range("K132").Value = 8833.1427
$ws.Range("M132").Value = -6303.1427

$ws.Range("H138").Value = 80406.164
$ws.Range("J138").Value = 80406.164
$ws.Range("L138").Value = 80406.164
$ws.Range("N138").Value = -90686.164

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8499
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 8499
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H74").Value = 12202.833
$ws.Range("J74").Value = 10312.25
$ws.Range("L74").Value = 10312.25
$ws.Range("N74").Value = -12184.25

$ws.Range("H77").Value = 12202.833
$ws.Range("J77").Value = 10312.25
$ws.Range("L77").Value = 30936.75
$ws.Range("N77").Value = -40296.75

$ws.Range("H100").Value = 55556344
$ws.Range("I100").Value = 557.8889
$ws.Range("K100").Value = 1115.7778
$ws.Range("M100").Value = -574.7778000000001

$ws.Range("H107").Value = 518.5
$ws.Range("I107").Value = 580
$ws.Range("J107").Value = 457
$ws.Range("K107").Value = 1740
$ws.Range("L107").Value = 1371
$ws.Range("M107").Value = 180
$ws.Range("N107").Value = -5211
